# Multiple Sim Verification Added
#
# 1. Update the test_files_path on the Main sheet (new dev machine).
# 2. Remove the obsolete "marker_area_file" row from every zDist_*
#    results sheet, shifting the rows below it up by one.
# 3. The "size50" and "size200" simulation sheets swap places (both
#    their tab names and their marker-size data row swap); the other
#    zDist_* sheets keep their identity.
# 4. Leave the "size1000" sheet active, matching the final save state.

$wb = $excel.ActiveWorkbook

# --- 1. Fix up the test files path on the Main sheet -----------------
$wsMain = $wb.Worksheets.Item("Main")
$wsMain.Range("B2").Value = "/home/paul/FiducialTags/Simulations"
[void]$wsMain.Range("B2").Select()

# --- 2. Drop the "marker_area_file" row (row 6) from every results ---
#        sheet; everything below shifts up automatically.
$zSheetNames = @(
    "zDist_size1000_camKhan",
    "zDist_size500_camKhan",
    "zDist_size50_camKhan",
    "zDist_size25_camKhan",
    "zDist_size200_camKhan",
    "zDist_size100_camKhan"
)
foreach ($name in $zSheetNames) {
    $ws = $wb.Worksheets.Item($name)
    [void]$ws.Rows.Item(6).Delete()
}

# --- 3. Swap the "size50" and "size200" sheets (name + marker file) --
$wsSize50 = $wb.Worksheets.Item("zDist_size50_camKhan")
$wsSize200 = $wb.Worksheets.Item("zDist_size200_camKhan")

$size50Marker = $wsSize50.Range("B7").Text
$size200Marker = $wsSize200.Range("B7").Text

$wsSize50.Range("B7").Value = $size200Marker
$wsSize200.Range("B7").Value = $size50Marker

$wsSize50.Name = "TEMP_SWAP_NAME_ZDIST"
$wsSize200.Name = "zDist_size50_camKhan"
$wsSize50.Name = "zDist_size200_camKhan"

# --- 4. Restore per-sheet selections left behind by the original -----
#        editing session.
[void]$wb.Worksheets.Item("zDist_size500_camKhan").Range("A8").Select()
[void]$wb.Worksheets.Item("zDist_size200_camKhan").Range("C29").Select()
[void]$wb.Worksheets.Item("zDist_size25_camKhan").Range("C77").Select()
[void]$wb.Worksheets.Item("zDist_size50_camKhan").Range("C17").Select()
[void]$wb.Worksheets.Item("zDist_size100_camKhan").Range("B7").Select()

# --- 5. Re-activate the "size1000" sheet ------------------------------
$wsSize1000 = $wb.Worksheets.Item("zDist_size1000_camKhan")
$wsSize1000.Activate()
[void]$wsSize1000.Range("B3").Select()
